$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45013, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Angeleno", "Especial", 280, 15000, 15000, 15000, "$/caja 15 kilos granel", "Paine", 1000, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45013, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Angeleno", "Primera", 300, 12000, 12000, 12000, "$/caja 15 kilos granel", "Paine", 800, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45013, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Angeleno", "Segunda", 290, 10000, 10000, 10000, "$/caja 15 kilos granel", "Paine", 667, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 210, 10000, 10000, 10000, "$/caja 15 kilos granel", "Región Metropolitana", 667, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Segunda", 280, 8000, 8000, 8000, "$/caja 15 kilos granel", "Región Metropolitana", 533, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Tercera", 170, 5000, 5000, 5000, "$/caja 15 kilos granel", "Región Metropolitana", 333, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Larry Ann", "Primera", 85, 10000, 10000, 10000, "$/caja 15 kilos granel", "Región de O'Higgins", 667, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44610, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Larry Ann", "Segunda", 120, 8000, 8000, 8000, "$/caja 15 kilos granel", "Región de O'Higgins", 533, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44568, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Especial", 280, 15000, 15000, 15000, "$/caja 15 kilos granel", "Región de O'Higgins", 1000, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44568, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 350, 12000, 12000, 12000, "$/caja 15 kilos granel", "Región de O'Higgins", 800, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44636, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Especial", 280, 11000, 11000, 11000, "$/caja 18 kilos granel", "Región de O'Higgins", 611, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44636, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 350, 9000, 9000, 9000, "$/caja 18 kilos granel", "Región de O'Higgins", 500, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44257, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Especial", 80, 11000, 11000, 11000, "$/caja 15 kilos granel", "Región de O'Higgins", 733, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44257, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 120, 8000, 8000, 8000, "$/caja 15 kilos granel", "Región de O'Higgins", 533, 15),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45008, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Angeleno", "Primera", 500, 12000, 12500, 12220, "$/caja 18 kilos granel", "Provincia de Curicó", 679, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44677, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Angeleno", "Primera", 220, 9000, 9000, 9000, "$/caja 18 kilos", "Región de O'Higgins", 500, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44677, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Angeleno", "Segunda", 250, 7200, 7200, 7200, "$/caja 18 kilos", "Región de O'Higgins", 400, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44209, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 70, 10000, 10000, 10000, "$/caja 18 kilos granel", "Región de O'Higgins", 556, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44209, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Segunda", 55, 8000, 8000, 8000, "$/caja 18 kilos granel", "Región de O'Higgins", 444, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44608, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Larry Ann", "Primera", 80, 12500, 12500, 12500, "$/bandeja 18 kilos granel", "Región Metropolitana", 694, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44608, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Larry Ann", "Segunda", 140, 10000, 10000, 10000, "$/bandeja 18 kilos granel", "Región Metropolitana", 556, 18),
  @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 44236, 13, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Black Amber", "Primera", 200, 9000, 10000, 9500, "$/caja 16 kilos granel", "Región de O'Higgins", 594, 16)
)

$startRow = 178
$r = $startRow
foreach ($row in $rows) {
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
    $r = $r + 1
}

# Apply date number format and style to column D for all affected rows (178-199)
$ws.Range("D178:D199").NumberFormat = $ws.Range("D2").NumberFormat

Write-Host "UsedRange:" $ws.UsedRange.Address()
